# Auto-generated edit script: updates specific market-data cells per the commit diff.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (index 1) ---
$ws = $wb.Worksheets.Item(1)
# Row 17
$ws.Cells.Item(17, 8).Value = 1820.1852  # H17
$ws.Cells.Item(17, 9).Value = 1064  # I17
$ws.Cells.Item(17, 10).Value = 2340.0625  # J17
$ws.Cells.Item(17, 11).Value = 3192  # K17
$ws.Cells.Item(17, 12).Value = 7020.1875  # L17
$ws.Cells.Item(17, 13).Value = -3024  # M17
$ws.Cells.Item(17, 14).Value = -7356.1875  # N17

# Row 62
$ws.Cells.Item(62, 8).Value = 621332.75  # H62
$ws.Cells.Item(62, 10).Value = 845832.5  # J62
$ws.Cells.Item(62, 12).Value = 845832.5  # L62
$ws.Cells.Item(62, 14).Value = -847080.5  # N62

# Row 65
$ws.Cells.Item(65, 8).Value = 621332.75  # H65
$ws.Cells.Item(65, 10).Value = 845832.5  # J65
$ws.Cells.Item(65, 12).Value = 4229162.5  # L65
$ws.Cells.Item(65, 14).Value = -4235402.5  # N65

# Row 109
$ws.Cells.Item(109, 8).Value = 89774  # H109
$ws.Cells.Item(109, 10).Value = 89774  # J109
$ws.Cells.Item(109, 12).Value = 89774  # L109
$ws.Cells.Item(109, 14).Value = -92548  # N109

# Row 110
$ws.Cells.Item(110, 8).Value = 67992  # H110
$ws.Cells.Item(110, 10).Value = 67992  # J110
$ws.Cells.Item(110, 12).Value = 67992  # L110
$ws.Cells.Item(110, 14).Value = -76172  # N110

# Row 117
$ws.Cells.Item(117, 8).Value = 90195.44500000001  # H117
$ws.Cells.Item(117, 10).Value = 90195.44500000001  # J117
$ws.Cells.Item(117, 12).Value = 90195.44500000001  # L117
$ws.Cells.Item(117, 14).Value = -99373.44500000001  # N117

# Row 120
$ws.Cells.Item(120, 8).Value = 49986.75  # H120
$ws.Cells.Item(120, 10).Value = 49986.75  # J120
$ws.Cells.Item(120, 12).Value = 49986.75  # L120
$ws.Cells.Item(120, 14).Value = -59662.75  # N120

# Row 123
$ws.Cells.Item(123, 8).Value = 76803.7  # H123
$ws.Cells.Item(123, 10).Value = 76803.7  # J123
$ws.Cells.Item(123, 12).Value = 76803.7  # L123
$ws.Cells.Item(123, 14).Value = -86603.7  # N123

# --- Sheet: ARM (index 2) ---
$ws = $wb.Worksheets.Item(2)
# Row 7
$ws.Cells.Item(7, 8).Value = 97307  # H7
$ws.Cells.Item(7, 10).Value = 97307  # J7
$ws.Cells.Item(7, 12).Value = 97307  # L7
$ws.Cells.Item(7, 14).Value = -97535  # N7

# Row 19
$ws.Cells.Item(19, 8).Value = 2452.5  # H19
$ws.Cells.Item(19, 9).Value = 3103.3333  # I19
$ws.Cells.Item(19, 11).Value = 3103.3333  # K19
$ws.Cells.Item(19, 13).Value = -2874.3333  # M19

# Row 26
$ws.Cells.Item(26, 13).Value = $null  # clear M26
$ws.Cells.Item(26, 8).Value = 0  # H26
$ws.Cells.Item(26, 9).Value = 0  # I26
$ws.Cells.Item(26, 11).Value = 0  # K26

# Row 32
$ws.Cells.Item(32, 8).Value = 10999.638  # H32
$ws.Cells.Item(32, 9).Value = 4961.846  # I32
$ws.Cells.Item(32, 11).Value = 4961.846  # K32
$ws.Cells.Item(32, 13).Value = -4674.846  # M32

# Row 45
$ws.Cells.Item(45, 8).Value = 2929.2307  # H45
$ws.Cells.Item(45, 9).Value = 2798.182  # I45
$ws.Cells.Item(45, 11).Value = 2798.182  # K45
$ws.Cells.Item(45, 13).Value = -2421.182  # M45

# Row 52
$ws.Cells.Item(52, 8).Value = 54686.6  # H52
$ws.Cells.Item(52, 10).Value = 54686.6  # J52
$ws.Cells.Item(52, 12).Value = 54686.6  # L52
$ws.Cells.Item(52, 14).Value = -55322.6  # N52

# Row 74
$ws.Cells.Item(74, 8).Value = 3112.4736  # H74
$ws.Cells.Item(74, 9).Value = 2539.5  # I74
$ws.Cells.Item(74, 10).Value = 3376.923  # J74
$ws.Cells.Item(74, 11).Value = 2539.5  # K74
$ws.Cells.Item(74, 12).Value = 3376.923  # L74
$ws.Cells.Item(74, 13).Value = -1665.5  # M74
$ws.Cells.Item(74, 14).Value = -5124.923  # N74

# Row 77
$ws.Cells.Item(77, 8).Value = 3112.4736  # H77
$ws.Cells.Item(77, 9).Value = 2539.5  # I77
$ws.Cells.Item(77, 10).Value = 3376.923  # J77
$ws.Cells.Item(77, 11).Value = 12697.5  # K77
$ws.Cells.Item(77, 12).Value = 16884.615  # L77
$ws.Cells.Item(77, 13).Value = -8329.5  # M77
$ws.Cells.Item(77, 14).Value = -25620.615  # N77

# Row 104
$ws.Cells.Item(104, 8).Value = 28961.4  # H104
$ws.Cells.Item(104, 10).Value = 28961.4  # J104
$ws.Cells.Item(104, 12).Value = 28961.4  # L104
$ws.Cells.Item(104, 14).Value = -35949.4  # N104

# Row 115
$ws.Cells.Item(115, 8).Value = 55491.332  # H115
$ws.Cells.Item(115, 10).Value = 55491.332  # J115
$ws.Cells.Item(115, 12).Value = 55491.332  # L115
$ws.Cells.Item(115, 14).Value = -58625.332  # N115

# Row 117
$ws.Cells.Item(117, 8).Value = 82686  # H117
$ws.Cells.Item(117, 10).Value = 82686  # J117
$ws.Cells.Item(117, 12).Value = 82686  # L117
$ws.Cells.Item(117, 14).Value = -91864  # N117

# Row 127
$ws.Cells.Item(127, 8).Value = 99991.8  # H127
$ws.Cells.Item(127, 10).Value = 99991.8  # J127
$ws.Cells.Item(127, 12).Value = 99991.8  # L127
$ws.Cells.Item(127, 14).Value = -109911.8  # N127

# --- Sheet: BSM (index 3) ---
$ws = $wb.Worksheets.Item(3)
# Row 13
$ws.Cells.Item(13, 8).Value = 43326.668  # H13
$ws.Cells.Item(13, 10).Value = 43326.668  # J13
$ws.Cells.Item(13, 12).Value = 43326.668  # L13
$ws.Cells.Item(13, 14).Value = -43662.668  # N13

# Row 51
$ws.Cells.Item(51, 8).Value = 74726.5  # H51
$ws.Cells.Item(51, 10).Value = 74726.5  # J51
$ws.Cells.Item(51, 12).Value = 74726.5  # L51
$ws.Cells.Item(51, 14).Value = -75708.5  # N51

# Row 53
$ws.Cells.Item(53, 8).Value = 49072.25  # H53
$ws.Cells.Item(53, 10).Value = 52526.668  # J53
$ws.Cells.Item(53, 12).Value = 52526.668  # L53
$ws.Cells.Item(53, 14).Value = -53674.668  # N53

# Row 108
$ws.Cells.Item(108, 8).Value = 99989  # H108
$ws.Cells.Item(108, 10).Value = 99989  # J108
$ws.Cells.Item(108, 12).Value = 99989  # L108
$ws.Cells.Item(108, 14).Value = -107669  # N108

# Row 109
$ws.Cells.Item(109, 8).Value = 99988  # H109
$ws.Cells.Item(109, 10).Value = 99988  # J109
$ws.Cells.Item(109, 12).Value = 99988  # L109
$ws.Cells.Item(109, 14).Value = -102762  # N109

# Row 114
$ws.Cells.Item(114, 8).Value = 89195.57000000001  # H114
$ws.Cells.Item(114, 10).Value = 88791.336  # J114
$ws.Cells.Item(114, 12).Value = 88791.336  # L114
$ws.Cells.Item(114, 14).Value = -97469.336  # N114

# Row 118
$ws.Cells.Item(118, 8).Value = 96352.664  # H118
$ws.Cells.Item(118, 10).Value = 95713.8  # J118
$ws.Cells.Item(118, 12).Value = 95713.8  # L118
$ws.Cells.Item(118, 14).Value = -99027.8  # N118

# Row 132
$ws.Cells.Item(132, 8).Value = 28024.838  # H132
$ws.Cells.Item(132, 10).Value = 28024.838  # J132
$ws.Cells.Item(132, 12).Value = 28024.838  # L132
$ws.Cells.Item(132, 14).Value = -38144.838  # N132

# Row 134
$ws.Cells.Item(134, 8).Value = 2360.5789  # H134
$ws.Cells.Item(134, 9).Value = 1665.6154  # I134
$ws.Cells.Item(134, 11).Value = 4996.8462  # K134
$ws.Cells.Item(134, 13).Value = -2461.8462  # M134

# --- Sheet: CRP (index 4) ---
$ws = $wb.Worksheets.Item(4)
# Row 9
$ws.Cells.Item(9, 8).Value = 37973.11  # H9
$ws.Cells.Item(9, 10).Value = 37973.11  # J9
$ws.Cells.Item(9, 12).Value = 37973.11  # L9
$ws.Cells.Item(9, 14).Value = -38309.11  # N9

# Row 16
$ws.Cells.Item(16, 8).Value = 3021.4285  # H16
$ws.Cells.Item(16, 9).Value = 2930  # I16
$ws.Cells.Item(16, 11).Value = 2930  # K16
$ws.Cells.Item(16, 13).Value = -2643  # M16

# Row 108
$ws.Cells.Item(108, 8).Value = 34726.125  # H108
$ws.Cells.Item(108, 10).Value = 34726.125  # J108
$ws.Cells.Item(108, 12).Value = 34726.125  # L108
$ws.Cells.Item(108, 14).Value = -42406.125  # N108

# Row 113
$ws.Cells.Item(113, 8).Value = 3021.4285  # H113
$ws.Cells.Item(113, 9).Value = 2930  # I113
$ws.Cells.Item(113, 11).Value = 2930  # K113
$ws.Cells.Item(113, 13).Value = -760  # M113

# Row 114
$ws.Cells.Item(114, 8).Value = 95847  # H114
$ws.Cells.Item(114, 10).Value = 95847  # J114
$ws.Cells.Item(114, 12).Value = 95847  # L114
$ws.Cells.Item(114, 14).Value = -104525  # N114

# Row 116
$ws.Cells.Item(116, 8).Value = 66380.8  # H116
$ws.Cells.Item(116, 10).Value = 66380.8  # J116
$ws.Cells.Item(116, 12).Value = 66380.8  # L116
$ws.Cells.Item(116, 14).Value = -75558.8  # N116

# Row 119
$ws.Cells.Item(119, 8).Value = 66966  # H119
$ws.Cells.Item(119, 10).Value = 66966  # J119
$ws.Cells.Item(119, 12).Value = 66966  # L119
$ws.Cells.Item(119, 14).Value = -76642  # N119

# Row 122
$ws.Cells.Item(122, 8).Value = 3038.4707  # H122
$ws.Cells.Item(122, 9).Value = 2708  # I122
$ws.Cells.Item(122, 11).Value = 8124  # K122
$ws.Cells.Item(122, 13).Value = -5674  # M122

# --- Sheet: CUL (index 5) ---
$ws = $wb.Worksheets.Item(5)
# Row 5
$ws.Cells.Item(5, 8).Value = 1461.9166  # H5
$ws.Cells.Item(5, 9).Value = 923.8333  # I5
$ws.Cells.Item(5, 10).Value = 2000  # J5
$ws.Cells.Item(5, 11).Value = 2771.4999  # K5
$ws.Cells.Item(5, 12).Value = 6000  # L5
$ws.Cells.Item(5, 13).Value = -2659.4999  # M5
$ws.Cells.Item(5, 14).Value = -6224  # N5

# Row 98
$ws.Cells.Item(98, 13).Value = $null  # clear M98
$ws.Cells.Item(98, 9).Value = 0  # I98
$ws.Cells.Item(98, 10).Value = 5224.5  # J98
$ws.Cells.Item(98, 11).Value = 0  # K98
$ws.Cells.Item(98, 12).Value = 15673.5  # L98
$ws.Cells.Item(98, 14).Value = -18669.5  # N98

# Row 135
$ws.Cells.Item(135, 8).Value = 1461.9166  # H135
$ws.Cells.Item(135, 9).Value = 923.8333  # I135
$ws.Cells.Item(135, 10).Value = 2000  # J135
$ws.Cells.Item(135, 11).Value = 8314.4997  # K135
$ws.Cells.Item(135, 12).Value = 18000  # L135
$ws.Cells.Item(135, 13).Value = -5779.4997  # M135
$ws.Cells.Item(135, 14).Value = -23070  # N135

# --- Sheet: GSM (index 6) ---
$ws = $wb.Worksheets.Item(6)
# Row 24
$ws.Cells.Item(24, 8).Value = 21166.666  # H24
$ws.Cells.Item(24, 10).Value = 17750  # J24
$ws.Cells.Item(24, 12).Value = 17750  # L24
$ws.Cells.Item(24, 14).Value = -18096  # N24

# Row 113
$ws.Cells.Item(113, 8).Value = 3829061.8  # H113
$ws.Cells.Item(113, 9).Value = 371137  # I113
$ws.Cells.Item(113, 11).Value = 371137  # K113
$ws.Cells.Item(113, 13).Value = -368967  # M113

# Row 114
$ws.Cells.Item(114, 8).Value = 65898.5  # H114
$ws.Cells.Item(114, 10).Value = 65898.5  # J114
$ws.Cells.Item(114, 12).Value = 65898.5  # L114
$ws.Cells.Item(114, 14).Value = -74576.5  # N114

# Row 119
$ws.Cells.Item(119, 8).Value = 69238.55  # H119
$ws.Cells.Item(119, 10).Value = 69293.3  # J119
$ws.Cells.Item(119, 12).Value = 69293.3  # L119
$ws.Cells.Item(119, 14).Value = -78969.3  # N119

# --- Sheet: LTW (index 7) ---
$ws = $wb.Worksheets.Item(7)
# Row 22
$ws.Cells.Item(22, 8).Value = 1342  # H22
$ws.Cells.Item(22, 9).Value = 967.5  # I22
$ws.Cells.Item(22, 10).Value = 1556  # J22
$ws.Cells.Item(22, 11).Value = 967.5  # K22
$ws.Cells.Item(22, 12).Value = 1556  # L22
$ws.Cells.Item(22, 13).Value = -672.5  # M22
$ws.Cells.Item(22, 14).Value = -2146  # N22

# Row 23
$ws.Cells.Item(23, 14).Value = $null  # clear N23
$ws.Cells.Item(23, 8).Value = 3000  # H23
$ws.Cells.Item(23, 9).Value = 3000  # I23
$ws.Cells.Item(23, 10).Value = 0  # J23
$ws.Cells.Item(23, 11).Value = 3000  # K23
$ws.Cells.Item(23, 12).Value = 0  # L23
$ws.Cells.Item(23, 13).Value = -2770  # M23

# Row 27
$ws.Cells.Item(27, 8).Value = 1342  # H27
$ws.Cells.Item(27, 9).Value = 967.5  # I27
$ws.Cells.Item(27, 10).Value = 1556  # J27
$ws.Cells.Item(27, 11).Value = 967.5  # K27
$ws.Cells.Item(27, 12).Value = 1556  # L27
$ws.Cells.Item(27, 13).Value = -860.5  # M27
$ws.Cells.Item(27, 14).Value = -1770  # N27

# Row 117
$ws.Cells.Item(117, 8).Value = 60480.4  # H117
$ws.Cells.Item(117, 10).Value = 60480.4  # J117
$ws.Cells.Item(117, 12).Value = 60480.4  # L117
$ws.Cells.Item(117, 14).Value = -69658.39999999999  # N117

# --- Sheet: WVR (index 8) ---
$ws = $wb.Worksheets.Item(8)
# Row 132
$ws.Cells.Item(132, 8).Value = 967361.5  # H132
$ws.Cells.Item(132, 9).Value = 943.6667  # I132
$ws.Cells.Item(132, 11).Value = 2831.0001  # K132
$ws.Cells.Item(132, 13).Value = -301.0001000000002  # M132
